$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2 through 12
# from serial date 45184 (2023-09-15) to 45185 (2023-09-16)
for ($row = 2; $row -le 12; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value2 = 45185
    }
}
